# Insert a new record row at row 65 on the worksheet, shifting the existing
# rows 65..146 down to 66..147. The new row captures a new Ají price record
# for "Americana (o)" in "Región del Maule".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65; this pushes everything below down.
$ws.Rows.Item(65).Insert()

# Populate the freshly inserted row 65 with the new data record.
$ws.Cells.Item(65, 1).Value = 5
$ws.Cells.Item(65, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(65, 3).Value = "Maule"
$ws.Cells.Item(65, 4).Value = 44546
$ws.Cells.Item(65, 5).Value = 7
$ws.Cells.Item(65, 6).Value = 100112021
$ws.Cells.Item(65, 7).Value = "Ají"
$ws.Cells.Item(65, 8).Value = "Americana (o)"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 100
$ws.Cells.Item(65, 11).Value = 17000
$ws.Cells.Item(65, 12).Value = 17000
$ws.Cells.Item(65, 13).Value = 17000
$ws.Cells.Item(65, 14).Value = "`$/caja 14 kilos"
$ws.Cells.Item(65, 15).Value = "Región del Maule"
$ws.Cells.Item(65, 16).Value = 1214
$ws.Cells.Item(65, 17).Value = 14
$ws.Cells.Item(65, 18).Value = "Hortaliza"
